# Commit message: "changing document, table attributes to lowerCamelCase"
#
# The only semantic (value) changes in the target diff are the two
# ObjTables header strings on the "!!Nodes" worksheet, where the
# attribute-like tokens are renamed from UpperCamelCase to lowerCamelCase:
#   A1: "!!!ObjTables ObjTablesVersion='0.0.8'" -> "!!!ObjTables objTablesVersion='0.0.8'"
#   A2: "!!ObjTables Type='Data' Id='Node'"      -> "!!ObjTables type='Data' id='Node'"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("!!Nodes")

$ws.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws.Range("A2").Value = "!!ObjTables type='Data' id='Node'"
